$d = $word.ActiveDocument

$d.Content.Find.Execute("95×32=", $true, $false, $false, $false, $false, $true, 1, $false, "99×98=", 2) | Out-Null
$d.Content.Find.Execute("20×23=", $true, $false, $false, $false, $false, $true, 1, $false, "43×22=", 2) | Out-Null
$d.Content.Find.Execute("83×27=", $true, $false, $false, $false, $false, $true, 1, $false, "67×89=", 2) | Out-Null
$d.Content.Find.Execute("42×95=", $true, $false, $false, $false, $false, $true, 1, $false, "51×28=", 2) | Out-Null
$d.Content.Find.Execute("59×92=", $true, $false, $false, $false, $false, $true, 1, $false, "18×14=", 2) | Out-Null
$d.Content.Find.Execute("25×92=", $true, $false, $false, $false, $false, $true, 1, $false, "47×41=", 2) | Out-Null
$d.Content.Find.Execute("83×17=", $true, $false, $false, $false, $false, $true, 1, $false, "60×66=", 2) | Out-Null
$d.Content.Find.Execute("38×14=", $true, $false, $false, $false, $false, $true, 1, $false, "69×89=", 2) | Out-Null
$d.Content.Find.Execute("82×99=", $true, $false, $false, $false, $false, $true, 1, $false, "17×39=", 2) | Out-Null
$d.Content.Find.Execute("81×77=", $true, $false, $false, $false, $false, $true, 1, $false, "91×65=", 2) | Out-Null
$d.Content.Find.Execute("49×14=", $true, $false, $false, $false, $false, $true, 1, $false, "41×18=", 2) | Out-Null
$d.Content.Find.Execute("63×67=", $true, $false, $false, $false, $false, $true, 1, $false, "35×64=", 2) | Out-Null
$d.Content.Find.Execute("54×33=", $true, $false, $false, $false, $false, $true, 1, $false, "19×87=", 2) | Out-Null
$d.Content.Find.Execute("67×99=", $true, $false, $false, $false, $false, $true, 1, $false, "58×39=", 2) | Out-Null
$d.Content.Find.Execute("60×38=", $true, $false, $false, $false, $false, $true, 1, $false, "21×65=", 2) | Out-Null
$d.Content.Find.Execute("95×35=", $true, $false, $false, $false, $false, $true, 1, $false, "43×29=", 2) | Out-Null
$d.Content.Find.Execute("45×39=", $true, $false, $false, $false, $false, $true, 1, $false, "70×79=", 2) | Out-Null
$d.Content.Find.Execute("85×60=", $true, $false, $false, $false, $false, $true, 1, $false, "73×61=", 2) | Out-Null
$d.Content.Find.Execute("79×75=", $true, $false, $false, $false, $false, $true, 1, $false, "40×82=", 2) | Out-Null
$d.Content.Find.Execute("17×55=", $true, $false, $false, $false, $false, $true, 1, $false, "92×11=", 2) | Out-Null
$d.Content.Find.Execute("13×83=", $true, $false, $false, $false, $false, $true, 1, $false, "96×86=", 2) | Out-Null
$d.Content.Find.Execute("51×49=", $true, $false, $false, $false, $false, $true, 1, $false, "70×59=", 2) | Out-Null
$d.Content.Find.Execute("22×26=", $true, $false, $false, $false, $false, $true, 1, $false, "38×62=", 2) | Out-Null
$d.Content.Find.Execute("57×95=", $true, $false, $false, $false, $false, $true, 1, $false, "83×36=", 2) | Out-Null
$d.Content.Find.Execute("26×32=", $true, $false, $false, $false, $false, $true, 1, $false, "42×78=", 2) | Out-Null
